$d = $word.ActiveDocument

# wdColorIndex constants
$wdYellow = 7
$wdRed = 6

function Set-ParaHighlight($para, $colorIndex) {
    # Going through Font (rather than Range directly) makes the
    # highlight fold into the paragraph-mark run properties (w:pPr/w:rPr)
    # as well as every run's properties (w:r/w:rPr) in the paragraph.
    $para.Range.Font.HighlightColorIndex = $colorIndex
}

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq "Логирование") {
        Set-ParaHighlight $p $wdYellow
    }
    elseif ($t -eq "Monitoring config") {
        Set-ParaHighlight $p $wdYellow
    }
    elseif ($t -eq "Hystrix") {
        Set-ParaHighlight $p $wdYellow
    }
    elseif ($t -eq "Hystrix – метрикс стрим") {
        Set-ParaHighlight $p $wdYellow
    }
    elseif ($t -eq "Hystrix - турбина") {
        Set-ParaHighlight $p $wdRed
    }
    elseif ($t -eq "ELK") {
        Set-ParaHighlight $p $wdRed
    }
}

Write-Output "done"
